# Weekly update: insert a new price record at row 93 for
# "Hortaliza, Feria Lagunitas de Puerto Montt - Haba", pushing the
# previously-existing rows 93..181 down to 94..182.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row before the current row 93 (shifts 93..181 -> 94..182,
# extends the used range from R181 to R182, and carries row 93's formatting
# down to the new row 94 - matching the style pattern needed for the date
# column in the new row 93 as well).
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new weekly record.
$ws.Cells.Item(93, 1).Value()  = 4
$ws.Cells.Item(93, 2).Value()  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(93, 3).Value()  = "Los Lagos"
$ws.Cells.Item(93, 4).Value()  = 45240
$ws.Cells.Item(93, 5).Value()  = 10
$ws.Cells.Item(93, 6).Value()  = 100112026
$ws.Cells.Item(93, 7).Value()  = "Haba"
$ws.Cells.Item(93, 8).Value()  = "Sin especificar"
$ws.Cells.Item(93, 9).Value()  = "Primera"
$ws.Cells.Item(93, 10).Value() = 180
$ws.Cells.Item(93, 11).Value() = 15000
$ws.Cells.Item(93, 12).Value() = 15000
$ws.Cells.Item(93, 13).Value() = 15000
$ws.Cells.Item(93, 14).Value() = "`$/saco 25 kilos"
$ws.Cells.Item(93, 15).Value() = "Región del Maule"
$ws.Cells.Item(93, 16).Value() = 600
$ws.Cells.Item(93, 17).Value() = 25
$ws.Cells.Item(93, 18).Value() = "Hortaliza"
